$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row of data (row 9) mirroring the structure of the existing rows
$ws.Range("A9").Value = "23-02-2019"
$ws.Range("B9").Value = 6
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 253
$ws.Range("F9").Value = 17
$ws.Range("G9").Value = 31
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 3
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 5
$ws.Range("M9").Value = 4
$ws.Range("N9").Value = 2
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = 12
$ws.Range("Q9").Value = 2
$ws.Range("R9").Value = 0
